$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.057.60"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.228.34"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.00%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.01%  "
$ws.Range("D15").Value = "2.568.69"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "2.222.21"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.734"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "39.990.46"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.19%  "
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.26%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("E35").Value = "  +2.55%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.10%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.112"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.78%  "
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").Value = "2.107.01"
$ws.Range("E42").Value = "  +9.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("E44").Value = "  +6.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").Value = "2.439.17"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  +6.52%  "
